$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.227.15'
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").Value = '2.780.99'
$ws.Range("E3").Value = '  +3.09%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '587.41'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '161.52'
$ws.Range("E6").Value = '  +8.56%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +2.37%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = '2.795.34'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").Value = '6.79'
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("D11").Value = '0.114'
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").Value = '0.400'
$ws.Range("E12").Value = '  +3.19%  '
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("D14").Value = '3.271.89'
$ws.Range("E14").Value = '  +2.94%  '
$ws.Range("D15").Value = '27.66'
$ws.Range("E15").Value = '  +4.39%  '
$ws.Range("D16").Value = '64.106.20'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '0.0000160'
$ws.Range("E17").Value = '  +6.14%  '
$ws.Range("D18").Value = '2.780.03'
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("D19").Value = '12.47'
$ws.Range("E19").Value = '  +4.72%  '
$ws.Range("D20").Value = '5.04'
$ws.Range("E20").Value = '  +3.52%  '
$ws.Range("D21").Value = '367.94'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("D22").Value = '7.08'
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("D23").Value = '0.576'
$ws.Range("E23").Value = '  +8.01%  '
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").Value = '67.62'
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("D26").Value = '0.177'
$ws.Range("E26").Value = '  +6.47%  '
$ws.Range("D27").Value = '8.78'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("D28").Value = '0.0₃0973'
$ws.Range("E28").Value = '  +14.37%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").Value = '7.32'
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("D32").Value = '1.27'
$ws.Range("E32").Value = '  +8.90%  '
$ws.Range("D33").Value = '172.95'
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '5.11'
$ws.Range("E34").Value = '  +7.39%  '
$ws.Range("D35").Value = '20.88'
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '1.49'
$ws.Range("E37").Value = '  +5.76%  '
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '1.03'
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").Value = '4.30'
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("D41").Value = '343.55'
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("D42").Value = '6.31'
$ws.Range("E42").Value = '  +10.71%  '
$ws.Range("D43").Value = '39.94'
$ws.Range("E43").Value = '  +2.49%  '
$ws.Range("D44").Value = '22.56'
$ws.Range("E44").Value = '  +4.70%  '
$ws.Range("D45").Value = '22.78'
$ws.Range("E45").Value = '  +4.57%  '
$ws.Range("D46").Value = '0.0613'
$ws.Range("E46").Value = '  +3.46%  '
$ws.Range("D47").Value = '0.654'
$ws.Range("E47").Value = '  +2.50%  '
$ws.Range("D48").Value = '0.0262'
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").Value = '138.96'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").Value = '0.103'
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("D51").Value = '2.178.73'
$ws.Range("E51").Value = '  +1.99%  '
